$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Write the 14 new backlog rows (19-32) as a single block write.
# ---------------------------------------------------------------------------
$arr = New-Object 'object[,]' 14,3
$arr[0,0] = "Instalar MySQL no servidor de dados da solução (VM Linux)"
$arr[0,1] = "Instalar o MySQL no servidor da máquina virtual para inserção dos dados."
$arr[0,2] = "Importante"
$arr[1,0] = "Inserção de dados do Arduino no MySQL (VM Linux)"
$arr[1,1] = "Inserção dos dados coletados pelos sensores para o banco de dados no servidor."
$arr[1,2] = "Importante"
$arr[2,0] = "Infraestrutura de cliente, servidor utilizando a VM Linux"
$arr[2,1] = "Infraestrutura do cliente funcionando com o servidor utilizando a máquina virtual."
$arr[2,2] = "Importante"
$arr[3,0] = "Planilha Sprint Backlog"
$arr[3,1] = "Criar uma planilha com um plano concreto de como atingir a Meta da Sprint, detalhando o trabalho a ser feito e como ele será realizado."
$arr[3,2] = "Importante"
$arr[4,0] = "Planilha Product Backlog"
$arr[4,1] = "Criar uma planilha com uma lista dinâmica e priorizada de todas as necessidades de um produto, como funcionalidades, melhorias e correções."
$arr[4,2] = "Importante"
$arr[5,0] = "Diagrama da Solução"
$arr[5,1] = "Criar um diagrama da solução com os principais elementos da solução, suas interações e como eles contribuem para atender às necessidades de negócio."
$arr[5,2] = "Importante"
$arr[6,0] = "Diagrama da Solução (Validado)"
$arr[6,1] = "Validar o diagrama de solução."
$arr[6,2] = "Importante"
$arr[7,0] = "Documento de Mudança"
$arr[7,1] = "Documento que registra cada mudança solicitada no decorrer do projeto e controla seu status."
$arr[7,2] = "Importante"
$arr[8,0] = "Ferramenta de Help Desk"
$arr[8,1] = "Criar e automatizar o processo de suporte e atendimento para solucionar problemas e responder a solicitações dos clientes."
$arr[8,2] = "Importante"
$arr[9,0] = "Fluxograma do suporte"
$arr[9,1] = "Criar o fluxograma de como funcionará o suporte ao cliente."
$arr[9,2] = "Importante"
$arr[10,0] = "Tela de Cadastro e Login"
$arr[10,1] = "Criar a tela de cadastro e login do site institucional."
$arr[10,2] = "Importante"
$arr[11,0] = "Tela de Dashboard"
$arr[11,1] = "Criar Tela de Dashboard dos dados coletados do site institucional."
$arr[11,2] = "Importante"
$arr[12,0] = "Tabelas criadas no MySQL - Protótipo - Final"
$arr[12,1] = "Protótipo das tabelas de banco de dados definidas para o banco de dados final."
$arr[12,2] = "Importante"
$arr[13,0] = "Tela de simulador financeiro - Final"
$arr[13,1] = "Tela de Simulador financeira definida para o site final."
$arr[13,2] = "Desejável"
$ws.Range("A19:C32").Value = $arr

# ---------------------------------------------------------------------------
# 2. Apply the same bordered / default-font formatting used by the existing
#    body rows (columns B/C of rows 3-18) to every new cell, by copying the
#    format from an already-formatted cell instead of constructing new
#    style entries from scratch.
# ---------------------------------------------------------------------------
$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("A19:C32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Column layout: unhide column B, widen both A and B (bestFit-like sizing)
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).Hidden = $false
$ws.Columns.Item(1).ColumnWidth = 50.17
$ws.Columns.Item(2).ColumnWidth = 193.17

# ---------------------------------------------------------------------------
# 4. View state: zoom to 85% and select the last populated cell (A31), mirroring
#    the author's final selection/scroll position when the sheet was saved.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("A31").Select() | Out-Null
